$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet (SCD0319 -> SCD0022)
$ws.Name = "SCD0022"

# 2. Update the TC_ID column (B) value for the data rows from "DGS-334" to "SCD0022-006"
$ws.Range("B2:B5").Value = "SCD0022-006"

# 3. Normalize alignment across the whole data table (left + vertically centered)
$table = $ws.Range("A1:R5")
$table.VerticalAlignment = -4108
$table.HorizontalAlignment = -4131

# 4. Column B needs to grow to fit the new, longer TC_ID text
$ws.Columns.Item(2).ColumnWidth = 11.5

# 5. Move the active selection/scroll position like in the final saved file
$ws.Range("B6").Select()
